$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-26 Sunday" "2025-10-27 Monday"

Replace-Text "56×32=1792" "48×64=3072"
Replace-Text "68×90=6120" "26×17=442"
Replace-Text "27×79=2133" "84×55=4620"
Replace-Text "85×95=8075" "47×92=4324"
Replace-Text "31×43=1333" "16×79=1264"

Replace-Text "68×40=2720" "98×98=9604"
Replace-Text "13×94=1222" "48×30=1440"
Replace-Text "99×87=8613" "31×22=682"
Replace-Text "43×22=946" "86×41=3526"
Replace-Text "71×16=1136" "23×86=1978"

Replace-Text "20×99=1980" "14×11=154"
Replace-Text "35×87=3045" "72×40=2880"
Replace-Text "43×56=2408" "38×66=2508"
Replace-Text "58×95=5510" "70×17=1190"
Replace-Text "16×42=672" "81×44=3564"

Replace-Text "30×63=1890" "94×84=7896"
Replace-Text "68×14=952" "76×47=3572"
Replace-Text "33×61=2013" "94×21=1974"
Replace-Text "29×13=377" "11×55=605"
Replace-Text "80×23=1840" "75×30=2250"

Replace-Text "89×30=2670" "23×77=1771"
Replace-Text "90×93=8370" "42×67=2814"
Replace-Text "30×59=1770" "90×33=2970"
Replace-Text "14×83=1162" "34×91=3094"
Replace-Text "16×74=1184" "96×28=2688"
